# Insert a new row at the top of the sheet, shifting all existing rows
# (including the former header row) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

# The inserted row 1 inherited no special formatting, while row 2 (the
# former header row) now carries the bold/centered/bordered header style.
# Move that formatting up to the new row 1 first...
$ws.Range("A2:L2").Copy() | Out-Null
$ws.Range("A1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ...then populate the new row 1 with the numeric column-index values.
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# Finally strip the header formatting back off row 2, which now holds the
# plain text column captions.
$ws.Range("A2:L2").ClearFormats()
